$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

# Refresh the per-gene query timestamps on the "data" sheet.
$data.Range("F2").Value = "2021-10-05 14:21:46.860602"
$data.Range("F3").Value = "2021-10-05 14:21:46.860610"
$data.Range("F4").Value = "2021-10-05 14:21:46.860613"
$data.Range("F5").Value = "2021-10-05 14:21:46.860616"
$data.Range("F6").Value = "2021-10-05 14:21:46.860618"
$data.Range("F7").Value = "2021-10-05 14:21:46.860621"
$data.Range("F8").Value = "2021-10-05 14:21:46.860624"
$data.Range("F9").Value = "2021-10-05 14:21:46.860626"
$data.Range("F10").Value = "2021-10-05 14:21:46.860629"
$data.Range("F11").Value = "2021-10-05 14:21:46.860632"
$data.Range("F12").Value = "2021-10-05 14:21:46.860634"
$data.Range("F13").Value = "2021-10-05 14:21:46.860637"

# Add a new "metadata" tab describing the panel this data was pulled from,
# placed after the "data" sheet.
$meta = $wb.Worksheets.Add()
$meta.Name = "metadata"
$meta.Move($null, $data)

$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Neurological segmental overgrowth"
$meta.Range("C2").Value = 736
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "1.15"
$meta.Range("E2").Value = "2021-03-26T16:54:03.247259Z"
$meta.Range("F2").Value = "2021-10-05 14:21:46.856935"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/736/?format=json"
